# Auto-generated edit script: updates cryptos list data cells (B/C/D/E) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.250.39"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.874.40"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'234.74"
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.4702"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").Value = "'0.2840"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.06610"
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").Value = "'20.26"
$ws.Range("E10").Value = "  +7.55%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'97.74"
$ws.Range("E12").Value = "  -4.11%  "
$ws.Range("D13").Value = "1.894.37"
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("D14").Value = "'5.086"
$ws.Range("E14").Value = "  -2.29%  "
$ws.Range("D15").Value = "'0.6727"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "'287.32"
$ws.Range("E16").Value = "  +7.45%  "
$ws.Range("D17").Value = "30.279.37"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.12%  "
$ws.Range("D19").Value = "'12.60"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "2.130.57"
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("D21").Value = "'5.386"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "'0.000007286"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("D23").Value = "'1.000"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("D25").Value = "'9.385"
$ws.Range("D26").Value = "'167.60"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'1.973"
$ws.Range("E28").Value = "  -4.42%  "
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("E31").Value = "  -5.22%  "
$ws.Range("D32").Value = "'1.466"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").Value = "'4.123"
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").Value = "'0.04693"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").Value = "'0.7064"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'0.9992"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "HuobiToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D38").Value = "'2.716"
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01872"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'6.446"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.525"
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'71.98"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.955"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8603"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'0.9996"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'103.01"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "TheSandbox"
$ws.Range("C47").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D47").Value = "'0.4196"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'990.16"
$ws.Range("E48").Value = "  +7.31%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.212"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.201"
$ws.Range("E50").Value = "  +4.37%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'34.05"
$ws.Range("E51").Value = "  -2.12%  "

# Reset style on force-texted numeric-looking cells so no stray style index is left behind
$numericRefs = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D14", "D15", "D16", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D28", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($ref in $numericRefs) {
    $ws.Range($ref).Style = "Normal"
}
